$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.379.60"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "1.837.61"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.56"
$ws.Range("E5").Value = "  -6.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5200"
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3219"
$ws.Range("E8").Value = "  -6.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06750"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.52"
$ws.Range("E10").Value = "  -8.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7619"
$ws.Range("E11").Value = "  -5.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07673"
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("D13").Value = "1.842.22"
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.52"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.016"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.93"
$ws.Range("E17").Value = "  -4.61%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007912"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "26.402.78"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "2.069.40"
$ws.Range("E21").Value = "  -3.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.550"
$ws.Range("E22").Value = "  -4.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.425"
$ws.Range("E23").Value = "  -6.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.925"
$ws.Range("E24").Value = "  -4.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.67"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.228"
$ws.Range("E26").Value = "  -5.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.653"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.165"
$ws.Range("E30").Value = "  -4.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.130"
$ws.Range("E31").Value = "  -4.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08729"
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04806"
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("E34").Value = "  -5.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.849"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6978"
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("E37").Value = "  -6.65%  "
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.187"
$ws.Range("E39").Value = "  -8.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4829"
$ws.Range("E40").Value = "  -6.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.15"
$ws.Range("E41").Value = "  -4.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8834"
$ws.Range("E42").Value = "  -7.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.085"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.613"
$ws.Range("E45").Value = "  -6.53%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05851"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4105"
$ws.Range("E47").Value = "  -8.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.982"
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("E49").Value = "  -4.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1215"
$ws.Range("E50").Value = "  -9.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8781"
$ws.Range("E51").Value = "  -0.74%  "

Write-Host "Applied cryptos update"
